$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 853.6667
$ws.Range("I32").Value = 400
$ws.Range("J32").Value = 910.375
$ws.Range("K32").Value = 400
$ws.Range("L32").Value = 910.375
$ws.Range("M32").Value = -74
$ws.Range("N32").Value = -1562.375

$ws.Range("H33").Value = 555.7727
$ws.Range("I33").Value = 116
$ws.Range("J33").Value = 1191
$ws.Range("K33").Value = 116
$ws.Range("L33").Value = 1191
$ws.Range("M33").Value = 113
$ws.Range("N33").Value = -1649

$ws.Range("H41").Value = 693.2857
$ws.Range("I41").Value = 783.6667
$ws.Range("J41").Value = 625.5
$ws.Range("K41").Value = 783.6667
$ws.Range("L41").Value = 625.5
$ws.Range("M41").Value = -343.6667
$ws.Range("N41").Value = -1505.5

$ws.Range("H58").Value = 2518.6667
$ws.Range("I58").Value = 111.42857
$ws.Range("J58").Value = 4625
$ws.Range("K58").Value = 334.28571
$ws.Range("L58").Value = 13875
$ws.Range("M58").Value = -184.28571
$ws.Range("N58").Value = -14175

$ws.Range("H74").Value = 3996.6667
$ws.Range("I74").Value = 3995
$ws.Range("K74").Value = 3995
$ws.Range("M74").Value = -3059

$ws.Range("H77").Value = 3996.6667
$ws.Range("I77").Value = 3995
$ws.Range("K77").Value = 19975
$ws.Range("M77").Value = -15295

$ws.Range("H113").Value = 3689.2307
$ws.Range("I113").Value = 3686
$ws.Range("J113").Value = 3700
$ws.Range("K113").Value = 3686
$ws.Range("L113").Value = 3700
$ws.Range("M113").Value = -432
$ws.Range("N113").Value = -10208

$ws.Range("H116").Value = 2133
$ws.Range("I116").Value = 1726.6666
$ws.Range("J116").Value = 2307.1428
$ws.Range("K116").Value = 1726.6666
$ws.Range("L116").Value = 2307.1428
$ws.Range("M116").Value = 1715.3334
$ws.Range("N116").Value = -9191.1428

$ws.Range("H121").Value = 1253.5
$ws.Range("J121").Value = 1253.5
$ws.Range("L121").Value = 3760.5
$ws.Range("N121").Value = -7254.5

$ws.Range("H137").Value = 3127680.8
$ws.Range("I137").Value = 3573392
$ws.Range("J137").Value = 7700.75
$ws.Range("K137").Value = 10720176
$ws.Range("L137").Value = 23102.25
$ws.Range("M137").Value = -10717626
$ws.Range("N137").Value = -28202.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2177.1667
$ws.Range("I2").Value = 2029.9231
$ws.Range("J2").Value = 2560
$ws.Range("K2").Value = 2029.9231
$ws.Range("L2").Value = 2560
$ws.Range("M2").Value = -1916.9231
$ws.Range("N2").Value = -2786

$ws.Range("H31").Value = 2471
$ws.Range("I31").Value = 2471
$ws.Range("K31").Value = 2471
$ws.Range("M31").Value = -2177

$ws.Range("H61").Value = 34552524
$ws.Range("I61").Value = 37074908
$ws.Range("K61").Value = 37074908
$ws.Range("M61").Value = -37074696

$ws.Range("H116").Value = 2177.1667
$ws.Range("I116").Value = 2029.9231
$ws.Range("J116").Value = 2560
$ws.Range("K116").Value = 2029.9231
$ws.Range("L116").Value = 2560
$ws.Range("M116").Value = 264.0769
$ws.Range("N116").Value = -7148

$ws.Range("H122").Value = 5293197.5
$ws.Range("I122").Value = 2264.8948
$ws.Range("J122").Value = 55557056
$ws.Range("K122").Value = 6794.6844
$ws.Range("L122").Value = 166671168
$ws.Range("M122").Value = -4344.6844
$ws.Range("N122").Value = -166676068

$ws.Range("H136").Value = 34552524
$ws.Range("I136").Value = 37074908
$ws.Range("K136").Value = 111224724
$ws.Range("M136").Value = -111222174

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2177.1667
$ws.Range("I3").Value = 2029.9231
$ws.Range("J3").Value = 2560
$ws.Range("K3").Value = 2029.9231
$ws.Range("L3").Value = 2560
$ws.Range("M3").Value = -1915.9231
$ws.Range("N3").Value = -2788

$ws.Range("H94").Value = 989.6667
$ws.Range("I94").Value = 986.7143
$ws.Range("K94").Value = 986.7143
$ws.Range("M94").Value = -535.7143

$ws.Range("H134").Value = 2536.7078
$ws.Range("I134").Value = 2229.8572
$ws.Range("J134").Value = 4446
$ws.Range("K134").Value = 6689.571599999999
$ws.Range("L134").Value = 13338
$ws.Range("M134").Value = -4154.571599999999
$ws.Range("N134").Value = -18408

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I58").Value = 62501596
$ws.Range("J58").Value = 3960.2
$ws.Range("K58").Value = 62501596
$ws.Range("L58").Value = 3960.2
$ws.Range("M58").Value = -62501393
$ws.Range("N58").Value = -4366.2

$ws.Range("H107").Value = 573.55554
$ws.Range("I107").Value = 520.25
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 520.25
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1399.75
$ws.Range("N107").Value = -4840

$ws.Range("H132").Value = 17937.732
$ws.Range("I132").Value = 1105.9608
$ws.Range("J132").Value = 113317.78
$ws.Range("K132").Value = 3317.8824
$ws.Range("L132").Value = 339953.34
$ws.Range("M132").Value = -787.8824000000004
$ws.Range("N132").Value = -345013.34

$ws.Range("H134").Value = 21740.768
$ws.Range("I134").Value = 1772.8372
$ws.Range("J134").Value = 87788.53999999999
$ws.Range("K134").Value = 5318.5116
$ws.Range("L134").Value = 263365.62
$ws.Range("M134").Value = -2783.5116
$ws.Range("N134").Value = -268435.62

$ws.Range("I136").Value = 62501596
$ws.Range("J136").Value = 3960.2
$ws.Range("K136").Value = 187504788
$ws.Range("L136").Value = 11880.6
$ws.Range("M136").Value = -187502238
$ws.Range("N136").Value = -16980.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I12").Value = 50000064
$ws.Range("J12").Value = 102.391304
$ws.Range("K12").Value = 150000192
$ws.Range("L12").Value = 307.173912
$ws.Range("M12").Value = -150000019
$ws.Range("N12").Value = -653.173912

$ws.Range("H131").Value = 1340
$ws.Range("J131").Value = 1370
$ws.Range("L131").Value = 4110
$ws.Range("N131").Value = -14190

$ws.Range("H136").Value = 1951.4445
$ws.Range("I136").Value = 1628.75
$ws.Range("J136").Value = 4533
$ws.Range("K136").Value = 4886.25
$ws.Range("L136").Value = 13599
$ws.Range("M136").Value = 213.75
$ws.Range("N136").Value = -23799

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1954.1666
$ws.Range("I113").Value = 1833.3334
$ws.Range("J113").Value = 2075
$ws.Range("K113").Value = 1833.3334
$ws.Range("L113").Value = 2075
$ws.Range("M113").Value = 336.6666
$ws.Range("N113").Value = -6415

$ws.Range("H132").Value = 72991.28999999999
$ws.Range("I132").Value = 41350.32
$ws.Range("K132").Value = 124050.96
$ws.Range("M132").Value = -121520.96

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 707.2857
$ws.Range("I22").Value = 700.1667
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 700.1667
$ws.Range("L22").Value = 750
$ws.Range("M22").Value = -405.1667
$ws.Range("N22").Value = -1340

$ws.Range("H27").Value = 707.2857
$ws.Range("I27").Value = 700.1667
$ws.Range("J27").Value = 750
$ws.Range("K27").Value = 700.1667
$ws.Range("L27").Value = 750
$ws.Range("M27").Value = -593.1667
$ws.Range("N27").Value = -964

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 30000
$ws.Range("J121").Value = 30000
$ws.Range("L121").Value = 30000
$ws.Range("N121").Value = -33494

$ws.Range("H136").Value = 38727.223
$ws.Range("I136").Value = 25087.38
$ws.Range("J136").Value = 86466.664
$ws.Range("K136").Value = 75262.14
$ws.Range("L136").Value = 259399.992
$ws.Range("M136").Value = -72712.14
$ws.Range("N136").Value = -264499.992
